$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 169430.6977861812
$ws.Range("C2").Value = 241893.9159574704
$ws.Range("D2").Value = 283060.0586359358
$ws.Range("E2").Value = 313191.8597465972

$ws.Range("B3").Value = 209135.9838772114
$ws.Range("C3").Value = 297219.0650641177
$ws.Range("D3").Value = 343277.4521467714
$ws.Range("E3").Value = 379730.7953214425

$ws.Range("B4").Value = 174098.9818665953
$ws.Range("C4").Value = 253160.8157205949
$ws.Range("D4").Value = 298968.5121198925
$ws.Range("E4").Value = 338343.8714689615

$ws.Range("B5").Value = 147464.2470118034
$ws.Range("C5").Value = 207532.1107858467
$ws.Range("D5").Value = 234189.4797150494
$ws.Range("E5").Value = 259128.2415958154

$ws.Range("B6").Value = 128733.317337992
$ws.Range("C6").Value = 180510.6698533363
$ws.Range("D6").Value = 205533.8476639566
$ws.Range("E6").Value = 225093.1384554107

$ws.Range("B7").Value = 14082.77606334089
$ws.Range("C7").Value = 19577.09992256755
$ws.Range("D7").Value = 22327.60211087606
$ws.Range("E7").Value = 24181.9562674005

$ws.Range("B8").Value = 685707.0930218007
$ws.Range("C8").Value = 978128.3897573499
$ws.Range("D8").Value = 1149091.596612659
$ws.Range("E8").Value = 1254415.841368718

$ws.Range("B9").Value = 194822.216171219
$ws.Range("C9").Value = 272649.7730737954
$ws.Range("D9").Value = 309371.3585756671
$ws.Range("E9").Value = 336415.5321210327

$ws.Range("B10").Value = 83634.44360114353
$ws.Range("C10").Value = 113168.5857715539
$ws.Range("D10").Value = 128230.0902425338
$ws.Range("E10").Value = 135030.7673417571

$ws.Range("B11").Value = 15240.88519867454
$ws.Range("C11").Value = 19544.27846766622
$ws.Range("D11").Value = 21889.30843569281
$ws.Range("E11").Value = 24883.75676046553

$ws.Range("B12").Value = 34115.22387689986
$ws.Range("C12").Value = 46293.74056691535
$ws.Range("D12").Value = 50055.43273309209
$ws.Range("E12").Value = 50610.91881612577

$ws.Range("B13").Value = 46880.94198278337
$ws.Range("C13").Value = 62802.34822430048
$ws.Range("D13").Value = 71774.72900780063
$ws.Range("E13").Value = 76450.36066295965
